$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value2 = 131054129
$ws.Range("B4").Value2 = 92267
$ws.Range("D4").Value2 = "VU"
$ws.Range("E4").Value2 = 1209
$ws.Range("F4").Value2 = "Rynkskinn"
$ws.Range("G4").Value2 = "Hermanssonia centrifuga"
$ws.Range("H4").Value2 = "(P. Karst.) Zmitr."
$ws.Range("P4").Value2 = "Rävmossen, Upl"
$ws.Range("Q4").Value2 = 660317
$ws.Range("R4").Value2 = 6661560
$ws.Range("S4").Value2 = 15
$ws.Range("T4").Value2 = "Uppsala"
$ws.Range("U4").Value2 = "Uppsala"
$ws.Range("V4").Value2 = "Uppland"
$ws.Range("W4").Value2 = "Rasbokil"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value2 = "2026-02-06"
$ws.Range("Z4").Value2 = "13:00"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value2 = "2026-02-06"
$ws.Range("AB4").Value2 = "13:00"
$ws.Range("AD4").Value2 = $false
$ws.Range("AE4").Value2 = $false
$ws.Range("AG4").Value2 = $false
$ws.Range("AW4").Value2 = "Tomas Falk"
$ws.Range("AX4").Value2 = "Tomas Falk"

# Row 5
$ws.Range("A5").Value2 = 131054159
$ws.Range("B5").Value2 = 92267
$ws.Range("D5").Value2 = "VU"
$ws.Range("E5").Value2 = 1209
$ws.Range("F5").Value2 = "Rynkskinn"
$ws.Range("G5").Value2 = "Hermanssonia centrifuga"
$ws.Range("H5").Value2 = "(P. Karst.) Zmitr."
$ws.Range("P5").Value2 = "Rävmossen, Upl"
$ws.Range("Q5").Value2 = 660304
$ws.Range("R5").Value2 = 6661793
$ws.Range("S5").Value2 = 15
$ws.Range("T5").Value2 = "Uppsala"
$ws.Range("U5").Value2 = "Uppsala"
$ws.Range("V5").Value2 = "Uppland"
$ws.Range("W5").Value2 = "Rasbokil"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value2 = "2026-02-06"
$ws.Range("Z5").Value2 = "10:18"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value2 = "2026-02-06"
$ws.Range("AB5").Value2 = "10:18"
$ws.Range("AD5").Value2 = $false
$ws.Range("AE5").Value2 = $false
$ws.Range("AG5").Value2 = $false
$ws.Range("AW5").Value2 = "Tomas Falk"
$ws.Range("AX5").Value2 = "Tomas Falk"

# Row 6
$ws.Range("A6").Value2 = 131054152
$ws.Range("B6").Value2 = 91808
$ws.Range("D6").Value2 = "NT"
$ws.Range("E6").Value2 = 1202
$ws.Range("F6").Value2 = "Ullticka"
$ws.Range("G6").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H6").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P6").Value2 = "Rävmossen, Upl"
$ws.Range("Q6").Value2 = 660473
$ws.Range("R6").Value2 = 6661613
$ws.Range("S6").Value2 = 15
$ws.Range("T6").Value2 = "Uppsala"
$ws.Range("U6").Value2 = "Uppsala"
$ws.Range("V6").Value2 = "Uppland"
$ws.Range("W6").Value2 = "Rasbokil"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value2 = "2026-02-06"
$ws.Range("Z6").Value2 = "10:53"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value2 = "2026-02-06"
$ws.Range("AB6").Value2 = "10:53"
$ws.Range("AC6").Value2 = "sumpskogsområde dikat"
$ws.Range("AD6").Value2 = $false
$ws.Range("AE6").Value2 = $false
$ws.Range("AG6").Value2 = $false
$ws.Range("AW6").Value2 = "Tomas Falk"
$ws.Range("AX6").Value2 = "Tomas Falk"

# Row 7
$ws.Range("A7").Value2 = 131054148
$ws.Range("B7").Value2 = 58043
$ws.Range("D7").Value2 = "NT"
$ws.Range("E7").Value2 = 103021
$ws.Range("F7").Value2 = "Talltita"
$ws.Range("G7").Value2 = "Poecile montanus"
$ws.Range("H7").Value2 = "(Conrad von Baldenstein, 1827)"
$ws.Range("P7").Value2 = "Rävmossen, Upl"
$ws.Range("Q7").Value2 = 660338
$ws.Range("R7").Value2 = 6661656
$ws.Range("S7").Value2 = 15
$ws.Range("T7").Value2 = "Uppsala"
$ws.Range("U7").Value2 = "Uppsala"
$ws.Range("V7").Value2 = "Uppland"
$ws.Range("W7").Value2 = "Rasbokil"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value2 = "2026-02-06"
$ws.Range("Z7").Value2 = "11:19"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value2 = "2026-02-06"
$ws.Range("AB7").Value2 = "11:19"
$ws.Range("AD7").Value2 = $false
$ws.Range("AE7").Value2 = $false
$ws.Range("AG7").Value2 = $false
$ws.Range("AW7").Value2 = "Tomas Falk"
$ws.Range("AX7").Value2 = "Tomas Falk"

# Row 8
$ws.Range("A8").Value2 = 131054160
$ws.Range("B8").Value2 = 91808
$ws.Range("D8").Value2 = "NT"
$ws.Range("E8").Value2 = 1202
$ws.Range("F8").Value2 = "Ullticka"
$ws.Range("G8").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H8").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P8").Value2 = "Rävmossen, Upl"
$ws.Range("Q8").Value2 = 660306
$ws.Range("R8").Value2 = 6661795
$ws.Range("S8").Value2 = 15
$ws.Range("T8").Value2 = "Uppsala"
$ws.Range("U8").Value2 = "Uppsala"
$ws.Range("V8").Value2 = "Uppland"
$ws.Range("W8").Value2 = "Rasbokil"
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value2 = "2026-02-06"
$ws.Range("Z8").Value2 = "10:18"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value2 = "2026-02-06"
$ws.Range("AB8").Value2 = "10:18"
$ws.Range("AD8").Value2 = $false
$ws.Range("AE8").Value2 = $false
$ws.Range("AG8").Value2 = $false
$ws.Range("AW8").Value2 = "Tomas Falk"
$ws.Range("AX8").Value2 = "Tomas Falk"

# Row 9
$ws.Range("A9").Value2 = 131054161
$ws.Range("B9").Value2 = 57881
$ws.Range("D9").Value2 = "NT"
$ws.Range("E9").Value2 = 100049
$ws.Range("F9").Value2 = "Spillkråka"
$ws.Range("G9").Value2 = "Dryocopus martius"
$ws.Range("H9").Value2 = "(Linnaeus, 1758)"
$ws.Range("P9").Value2 = "Rävmossen, Upl"
$ws.Range("Q9").Value2 = 660212
$ws.Range("R9").Value2 = 6661786
$ws.Range("S9").Value2 = 15
$ws.Range("T9").Value2 = "Uppsala"
$ws.Range("U9").Value2 = "Uppsala"
$ws.Range("V9").Value2 = "Uppland"
$ws.Range("W9").Value2 = "Rasbokil"
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value2 = "2026-02-06"
$ws.Range("Z9").Value2 = "10:05"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value2 = "2026-02-06"
$ws.Range("AB9").Value2 = "10:05"
$ws.Range("AD9").Value2 = $false
$ws.Range("AE9").Value2 = $false
$ws.Range("AG9").Value2 = $false
$ws.Range("AW9").Value2 = "Tomas Falk"
$ws.Range("AX9").Value2 = "Tomas Falk"

# Row 10
$ws.Range("A10").Value2 = 131054128
$ws.Range("B10").Value2 = 4779
$ws.Range("D10").Value2 = "LC"
$ws.Range("E10").Value2 = 102306
$ws.Range("F10").Value2 = "Granbarkgnagare"
$ws.Range("G10").Value2 = "Microbregma emarginatum"
$ws.Range("H10").Value2 = "(Duftschmid, 1825)"
$ws.Range("P10").Value2 = "Rävmossen, Upl"
$ws.Range("Q10").Value2 = 660278
$ws.Range("R10").Value2 = 6661596
$ws.Range("S10").Value2 = 15
$ws.Range("T10").Value2 = "Uppsala"
$ws.Range("U10").Value2 = "Uppsala"
$ws.Range("V10").Value2 = "Uppland"
$ws.Range("W10").Value2 = "Rasbokil"
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value2 = "2026-02-06"
$ws.Range("Z10").Value2 = "13:17"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value2 = "2026-02-06"
$ws.Range("AB10").Value2 = "13:17"
$ws.Range("AD10").Value2 = $false
$ws.Range("AE10").Value2 = $false
$ws.Range("AG10").Value2 = $false
$ws.Range("AW10").Value2 = "Tomas Falk"
$ws.Range("AX10").Value2 = "Tomas Falk"

# Row 11
$ws.Range("A11").Value2 = 131054149
$ws.Range("B11").Value2 = 91808
$ws.Range("D11").Value2 = "NT"
$ws.Range("E11").Value2 = 1202
$ws.Range("F11").Value2 = "Ullticka"
$ws.Range("G11").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H11").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P11").Value2 = "Rävmossen, Upl"
$ws.Range("Q11").Value2 = 660399
$ws.Range("R11").Value2 = 6661680
$ws.Range("S11").Value2 = 15
$ws.Range("T11").Value2 = "Uppsala"
$ws.Range("U11").Value2 = "Uppsala"
$ws.Range("V11").Value2 = "Uppland"
$ws.Range("W11").Value2 = "Rasbokil"
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value2 = "2026-02-06"
$ws.Range("Z11").Value2 = "11:01"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value2 = "2026-02-06"
$ws.Range("AB11").Value2 = "11:01"
$ws.Range("AD11").Value2 = $false
$ws.Range("AE11").Value2 = $false
$ws.Range("AG11").Value2 = $false
$ws.Range("AW11").Value2 = "Tomas Falk"
$ws.Range("AX11").Value2 = "Tomas Falk"
